$wb = $excel.ActiveWorkbook

$oldText = "February 03 2026 17.29.55 EST"
$newText = "February 03 2026 18.05.36 EST"

$wsAbout = $wb.Worksheets.Item("About")

$cellA2 = $wsAbout.Range("A2")
$a2 = $cellA2.Value()
$cellA2.Value = $a2.Replace($oldText, $newText)

$cellA6 = $wsAbout.Range("A6")
$a6 = $cellA6.Value()
$cellA6.Value = $a6.Replace($oldText, $newText)

$wsData = $wb.Worksheets.Item("Boundaries and methane sources")
for ($row = 2; $row -le 8; $row++) {
    $cell = $wsData.Cells.Item($row, 19)  # column S
    $val = $cell.Value()
    $cell.Value = $val.Replace($oldText, $newText)
}
